$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 24, shifting existing rows 24-40 down to 25-41.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with data (a new weekly price entry),
# mirroring the columns of the former row 24 (now row 25) except for the
# updated Date / Volumen / Precio minimo / Precio maximo / Precio promedio
# ponderado / Precio $/Kg values.
$ws.Range("A24").Value = 9
$ws.Range("B24").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C24").Value = "Metropolitana"
$ws.Range("D24").Value = 44830
$ws.Range("E24").Value = 13
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100108
$ws.Range("H24").Value = "Tropicales y subtropicales"
$ws.Range("I24").Value = 100108007
$ws.Range("J24").Value = "Coco"
$ws.Range("K24").Value = "Sin especificar"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 200
$ws.Range("N24").Value = 30000
$ws.Range("O24").Value = 30000
$ws.Range("P24").Value = 30000
$ws.Range("Q24").Value = "$/malla 20 unidades"
$ws.Range("R24").Value = "Perú"
$ws.Range("S24").Value = 1500
$ws.Range("T24").Value = 20
